$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column retains exact text formatting (avoid Excel auto-numeric conversion)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.150.02"
$ws.Range("E2").Value = "  -1.96%  "
$ws.Range("D3").Value = "3.236.88"
$ws.Range("E3").Value = "  -5.25%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "581.63"
$ws.Range("E5").Value = "  -4.76%  "
$ws.Range("D6").Value = "143.21"
$ws.Range("E6").Value = "  -13.70%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.223.97"
$ws.Range("E8").Value = "  -5.42%  "
$ws.Range("D9").Value = "0.528"
$ws.Range("E9").Value = "  -11.13%  "
$ws.Range("D10").Value = "0.164"
$ws.Range("E10").Value = "  -14.62%  "
$ws.Range("D11").Value = "6.67"
$ws.Range("E11").Value = "  -2.76%  "
$ws.Range("D12").Value = "0.487"
$ws.Range("E12").Value = "  -12.62%  "
$ws.Range("D13").Value = "0.0000237"
$ws.Range("E13").Value = "  -10.98%  "
$ws.Range("D14").Value = "36.74"
$ws.Range("E14").Value = "  -15.91%  "
$ws.Range("D15").Value = "3.749.61"
$ws.Range("E15").Value = "  -5.50%  "
$ws.Range("D16").Value = "67.168.28"
$ws.Range("E16").Value = "  -2.21%  "
$ws.Range("D17").Value = "3.237.92"
$ws.Range("E17").Value = "  -5.33%  "
$ws.Range("D18").Value = "0.112"
$ws.Range("E18").Value = "  -6.41%  "
$ws.Range("D19").Value = "6.87"
$ws.Range("E19").Value = "  -14.72%  "
$ws.Range("D20").Value = "505.50"
$ws.Range("E20").Value = "  -12.02%  "
$ws.Range("D21").Value = "14.49"
$ws.Range("E21").Value = "  -14.56%  "
$ws.Range("D22").Value = "0.730"
$ws.Range("E22").Value = "  -12.94%  "
$ws.Range("D23").Value = "7.48"
$ws.Range("E23").Value = "  -16.20%  "
$ws.Range("D24").Value = "82.99"
$ws.Range("E24").Value = "  -12.20%  "
$ws.Range("D25").Value = "12.98"
$ws.Range("E25").Value = "  -13.19%  "
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D27").Value = "3.12"
$ws.Range("E27").Value = "  -13.22%  "
$ws.Range("D28").Value = "2.07"
$ws.Range("E28").Value = "  -13.04%  "
$ws.Range("D29").Value = "7.70"
$ws.Range("E29").Value = "  -9.07%  "
$ws.Range("D30").Value = "28.07"
$ws.Range("E30").Value = "  -13.29%  "
$ws.Range("D31").Value = "1.16"
$ws.Range("E31").Value = "  -5.06%  "
$ws.Range("D32").Value = "2.54"
$ws.Range("E32").Value = "  -7.58%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  -0.20%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "6.25"
$ws.Range("E34").Value = "  -19.34%  "
$ws.Range("D35").Value = "5.45"
$ws.Range("E35").Value = "  -16.13%  "
$ws.Range("D36").Value = "54.93"
$ws.Range("E36").Value = "  -1.87%  "
$ws.Range("D37").Value = "493.37"
$ws.Range("E37").Value = "  -15.29%  "
$ws.Range("D38").Value = "0.0425"
$ws.Range("E38").Value = "  -7.91%  "
$ws.Range("D39").Value = "0.0823"
$ws.Range("E39").Value = "  -12.99%  "
$ws.Range("D40").Value = "8.64"
$ws.Range("E40").Value = "  -16.81%  "
$ws.Range("E41").Value = "  -12.99%  "
$ws.Range("D42").Value = "2.873.49"
$ws.Range("E42").Value = "  -10.38%  "
$ws.Range("D43").Value = "2.58"
$ws.Range("E43").Value = "  -14.49%  "
$ws.Range("D44").Value = "0.256"
$ws.Range("E44").Value = "  -11.85%  "
$ws.Range("D46").Value = "2.10"
$ws.Range("E46").Value = "  -10.64%  "
$ws.Range("D47").Value = "25.63"
$ws.Range("E47").Value = "  -17.44%  "
$ws.Range("D48").Value = "0.0₃0537"
$ws.Range("E48").Value = "  -19.27%  "
$ws.Range("D49").Value = "122.26"
$ws.Range("E49").Value = "  -7.62%  "
$ws.Range("D50").Value = "0.110"
$ws.Range("E50").Value = "  -12.00%  "
$ws.Range("D51").Value = "2.21"
$ws.Range("E51").Value = "  -19.81%  "
